# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" message text ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$lines = @(
    "Conversión del día 💰",
    "✅ Dólar paralelo: 68",
    "",
    "Binance",
    "✅ 1000 Bs = 8.7 = 36227.66 pesos",
    "✅ 36227.66 pesos = 8.68 = 955.19 Bs",
    "",
    "Promedio competencia",
    "✅ Tasa pesos: 20",
    "✅ Tasa Bs: 20",
    "✅ % Ganancia: 20%"
)
$newText = [string]::Join("`n", $lines)

$ws1.Range("A1").Value = $newText

# --- tasas: update N10, O10 and N12 values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 114.995
$ws2.Range("O10").Value = 4166
$ws2.Range("N12").Value = 4172
